$wb = $excel.ActiveWorkbook

# Data rows: SheetIndex,RowNumber,ColumnNumber,NewValue ("NULL" = clear/remove the cell)
$changesText = @"
1,4,8,400.66666
1,4,9,301.75
1,4,11,301.75
1,4,13,-187.75
1,51,8,10298.521
1,51,9,8985.714
1,51,10,10872.875
1,51,11,8985.714
1,51,12,10872.875
1,51,13,-8501.714
1,51,14,-11840.875
1,53,8,7804.3125
1,53,9,761.7778
1,53,11,761.7778
1,53,13,-124.7778
1,69,8,7653.25
1,69,10,8833.333000000001
1,69,12,26499.999
1,69,14,-28247.999
1,72,8,7653.25
1,72,10,8833.333000000001
1,72,12,79499.997
1,72,14,-88235.997
1,118,8,523.2222
1,118,9,234.83333
1,118,10,1100
1,118,11,704.49999
1,118,12,3300
1,118,13,952.50001
1,118,14,-6614
1,125,8,5194.9165
1,125,9,3849
1,125,10,20000
1,125,11,34641
1,125,12,180000
1,125,13,-32181
1,125,14,-184920
1,127,8,1000896.8
1,127,9,1000896.8
1,127,11,3002690.4
1,127,13,-2997730.4
1,141,8,3306.75
1,141,9,2411.0667
1,141,10,5993.8
1,141,11,7233.2001
1,141,12,17981.4
1,141,13,-2053.2001
1,141,14,-28341.4
2,8,8,2510620
2,8,9,3334160
2,8,11,3334160
2,8,13,-3334016
2,10,8,4
2,10,9,4
2,10,11,4
2,10,13,166
2,12,8,1836.3636
2,12,9,1836.3636
2,12,11,1836.3636
2,12,13,-1663.3636
2,13,8,251083.83
2,13,9,375000.75
2,13,11,375000.75
2,13,13,-374856.75
2,61,8,15314.143
2,61,9,3064.2222
2,61,11,3064.2222
2,61,13,-2852.2222
2,74,8,20134.912
2,74,9,2814.1
2,74,10,33458.617
2,74,11,2814.1
2,74,12,33458.617
2,74,13,-1940.1
2,74,14,-35206.617
2,77,8,20134.912
2,77,9,2814.1
2,77,10,33458.617
2,77,11,14070.5
2,77,12,167293.085
2,77,13,-9702.5
2,77,14,-176029.085
2,132,8,2640812.8
2,132,9,3285.2903
2,132,11,9855.8709
2,132,13,-7325.8709
2,135,8,129499.5
2,135,10,129499.5
2,135,12,129499.5
2,135,14,-139639.5
2,136,8,15314.143
2,136,9,3064.2222
2,136,11,9192.6666
2,136,13,-6642.6666
3,25,8,1785.6666
3,25,10,9000
3,25,12,9000
3,25,14,-9470
3,40,8,49999.75
3,40,10,49999.75
3,40,12,49999.75
3,40,14,-50529.75
3,74,8,0
3,74,10,0
3,74,12,0
3,74,14,NULL
3,77,8,0
3,77,10,0
3,77,12,0
3,77,14,NULL
3,99,8,1399.3
3,99,9,1097.5
3,99,11,1097.5
3,99,13,400.5
3,134,8,13136.029
3,134,9,7731
3,134,10,19216.688
3,134,11,23193
3,134,12,57650.064
3,134,13,-20658
3,134,14,-62720.064
4,8,8,1008.5
4,8,10,1338.3334
4,8,12,1338.3334
4,8,14,-1618.3334
4,10,8,91080
4,10,9,170.875
4,10,10,333504.34
4,10,11,170.875
4,10,12,333504.34
4,10,13,-31.875
4,10,14,-333782.34
4,11,8,266.33334
4,11,10,266.33334
4,11,12,266.33334
4,11,14,-546.33334
4,12,8,1662.8235
4,12,9,944.3333
4,12,10,2471.125
4,12,11,944.3333
4,12,12,2471.125
4,12,13,-774.3333
4,12,14,-2811.125
4,13,8,2000
4,13,10,2000
4,13,12,2000
4,13,14,-2278
4,17,8,4336
4,17,9,8
4,17,10,6500
4,17,11,8
4,17,12,6500
4,17,13,166
4,17,14,-6848
4,31,8,20189.13
4,31,9,6203.2383
4,31,10,37465.824
4,31,11,6203.2383
4,31,12,37465.824
4,31,13,-5908.2383
4,31,14,-38055.824
4,34,8,20189.13
4,34,9,6203.2383
4,34,10,37465.824
4,34,11,6203.2383
4,34,12,37465.824
4,34,13,-6001.2383
4,34,14,-37869.824
4,58,8,11927.766
4,58,9,4122.2383
4,58,10,18232.23
4,58,11,4122.2383
4,58,12,18232.23
4,58,13,-3919.2383
4,58,14,-18638.23
4,122,8,5115.8076
4,122,9,2215.2354
4,122,10,10594.667
4,122,11,6645.706200000001
4,122,12,31784.001
4,122,13,-4195.706200000001
4,122,14,-36684.001
4,132,8,6526.407
4,132,9,2065.647
4,132,10,14109.7
4,132,11,6196.941
4,132,12,42329.10000000001
4,132,13,-3666.941
4,132,14,-47389.10000000001
4,134,8,29417932
4,134,9,1739.6842
4,134,10,66678444
4,134,11,5219.0526
4,134,12,200035332
4,134,13,-2684.0526
4,134,14,-200040402
4,136,8,11927.766
4,136,9,4122.2383
4,136,10,18232.23
4,136,11,12366.7149
4,136,12,54696.69
4,136,13,-9816.714899999999
4,136,14,-59796.69
5,4,8,93459090
5,4,9,128380820
5,4,11,385142460
5,4,13,-385142348
5,11,8,365.375
5,11,9,137.33333
5,11,11,411.99999
5,11,13,-271.99999
5,39,8,3499.9092
5,39,10,10000
5,39,12,30000
5,39,14,-30588
5,55,8,3960.4375
5,55,10,3697.4167
5,55,12,11092.2501
5,55,14,-11446.2501
5,80,8,18580.133
5,80,9,1400.5
5,80,11,4201.5
5,80,13,-3265.5
5,83,8,18580.133
5,83,9,1400.5
5,83,11,12604.5
5,83,13,-7924.5
5,103,8,6795.8184
5,103,9,350.8
5,103,10,12166.667
5,103,11,1052.4
5,103,12,36500.001
5,103,13,-173.4000000000001
5,103,14,-38258.001
5,108,8,666
5,108,9,666
5,108,11,1998
5,108,13,882
5,114,8,1688.3636
5,114,9,1366.25
5,114,10,1872.4286
5,114,11,4098.75
5,114,12,5617.2858
5,114,13,-844.75
5,114,14,-12125.2858
5,117,8,2227
5,117,9,175
5,117,10,2455
5,117,11,525
5,117,12,7365
5,117,13,2917
5,117,14,-14249
5,131,8,1491.63
5,131,10,1491.63
5,131,12,4474.89
5,131,14,-14554.89
6,9,8,236.77777
6,9,9,105.666664
6,9,11,105.666664
6,9,13,64.333336
6,10,8,5000
6,10,10,5000
6,10,12,5000
6,10,14,-5338
6,11,8,650902.2
6,11,10,6220.1113
6,11,12,6220.1113
6,11,14,-6498.1113
6,12,8,775
6,12,10,366.66666
6,12,12,366.66666
6,12,14,-646.66666
6,49,8,37815.5
6,49,10,37815.5
6,49,12,37815.5
6,49,14,-38183.5
6,69,8,59122.5
6,69,10,59122.5
6,69,12,59122.5
6,69,14,-60620.5
6,72,8,59122.5
6,72,10,59122.5
6,72,12,177367.5
6,72,14,-184855.5
6,132,8,11725.615
6,132,9,12440.292
6,132,11,37320.876
6,132,13,-34790.876
7,7,8,10484
7,7,9,2196.5
7,7,10,11990.818
7,7,11,2196.5
7,7,12,11990.818
7,7,13,-2084.5
7,7,14,-12214.818
7,21,8,0
7,21,9,0
7,21,11,0
7,21,13,NULL
7,25,8,1949998
7,25,10,2999998
7,25,12,2999998
7,25,14,-3000458
7,40,8,20399.8
7,40,9,7000
7,40,10,23749.75
7,40,11,7000
7,40,12,23749.75
7,40,13,-6864
7,40,14,-24021.75
7,55,8,2180.5
7,55,9,1399
7,55,11,1399
7,55,13,-1226
7,61,8,4091.7273
7,61,9,2659.077
7,61,11,2659.077
7,61,13,-2457.077
7,113,8,4091.7273
7,113,9,2659.077
7,113,11,2659.077
7,113,13,-489.0770000000002
7,122,8,8305.956
7,122,9,5964.5386
7,122,11,17893.6158
7,122,13,-15443.6158
7,126,8,10484
7,126,9,2196.5
7,126,10,11990.818
7,126,11,6589.5
7,126,12,35972.454
7,126,13,-4119.5
7,126,14,-40912.454
7,132,8,897655
7,132,9,3598.45
7,132,11,10795.35
7,132,13,-8265.349999999999
8,34,8,9999.333000000001
8,34,10,9999
8,34,12,9999
8,34,14,-10405
8,113,8,3256.0356
8,113,9,3513.7778
8,113,10,2792.1
8,113,11,10541.3334
8,113,12,8376.299999999999
8,113,13,-8371.3334
8,113,14,-12716.3
8,122,8,9368.467000000001
8,122,9,6058.5
8,122,10,10572.091
8,122,11,18175.5
8,122,12,31716.273
8,122,13,-15725.5
8,122,14,-36616.273
"@

$changes = $changesText -split "`n"
$wsCache = @{}
$count = 0

foreach ($line in $changes) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line.Split(",")
    $sheetIdx = [int]$parts[0]
    $row = [int]$parts[1]
    $col = [int]$parts[2]
    $valStr = $parts[3]

    if (-not $wsCache.ContainsKey($sheetIdx)) {
        $wsCache[$sheetIdx] = $wb.Worksheets.Item($sheetIdx)
    }
    $ws = $wsCache[$sheetIdx]

    if ($valStr -eq "NULL") {
        $ws.Cells.Item($row, $col).Value = $null
    } else {
        $v = $valStr -as [double]
        $ws.Cells.Item($row, $col).Value = $v
    }
    $count = $count + 1
}

Write-Host "Applied $count cell changes"
